# Updated symbol list on Tue Dec 27 09:22:02 UTC 2022 with GitHub Actions
#
# Refresh the crypto price table: column D ("Price") gets new quotes for a
# number of rows, a couple of "Worstin24h" suffix flags move between the
# Volume(1h) labels (E22 / E43), and rows 41-43 get re-ranked (KickToken,
# BKEXToken and CEJI shift rank/price/link).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") updates -------------------------------------------
# These are stored as plain text in the sheet (not numbers), and some of the
# new values have significant trailing zeros (e.g. "3.430", "0.0002000") or
# very small magnitudes that Excel would otherwise mangle/round if it were
# allowed to auto-coerce them to a numeric type. Forcing the cell to Text
# format before assigning the value keeps the exact string representation.
$priceUpdates = @{
    "D2"  = "242.74"
    "D4"  = "5.406"
    "D5"  = "0.05947"
    "D6"  = "3.430"
    "D7"  = "6.498"
    "D8"  = "0.8146"
    "D9"  = "0.9186"
    "D10" = "0.1435"
    "D11" = "0.07412"
    "D12" = "0.03281"
    "D13" = "0.03067"
    "D14" = "0.09348"
    "D15" = "3.849"
    "D16" = "0.001588"
    "D17" = "0.04687"
    "D18" = "0.0005951"
    "D19" = "0.005903"
    "D20" = "0.001263"
    "D21" = "0.004793"
    "D22" = "0.00008003"
    "D23" = "3.576"
    "D25" = "0.3239"
    "D40" = "0.03938"
    "D41" = "0.1073"
    "D42" = "0.002551"
    "D43" = "0.003068"
    "D44" = "0.008906"
    "D45" = "0.00005171"
    "D47" = "0.7192"
    "D49" = "0.00002100"
    "D50" = "0.0002000"
}
foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = "Normal"
}

# --- Row 22: NitroEx no longer flagged as the worst performer --------------
$ws.Range("E22").Value = "21NitroExNTX"

# --- Rows 41-43: coins re-ranked --------------------------------------------
# Row 41 becomes BKEXToken (was KickToken)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42 becomes CEJI (was BKEXToken)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 becomes KickToken (was CEJI), and picks up the "Worstin24h" flag
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
